$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.229.48'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '2.935.12'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +4.02%  '
$ws.Range("E10").Value = '  -0.78%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '3.421.39'
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").Value = '61.192.91'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '2.938.07'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '433.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -3.48%  '
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("E31").Value = '  +3.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '0.0₃0879'
$ws.Range("E34").Value = '  +2.96%  '
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.64'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '41.36'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.284'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '377.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0348'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '2.706.10'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("E50").Value = '  -0.40%  '
$ws.Range("E51").Value = '  +0.42%  '
